# Arbeitsjournal KW07 -> KW09 edits
# (wording clean-up: remove "(Schnupperer)" gloss next to "Nico", merge a
#  few runs that had no real formatting differences, and rework the
#  Wochenrueckblick paragraph so it names "Nico" instead of "ein Schnupperer")
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Table cell: "Auf Nico (Schnupperer) aufgepasst und bewertet"
#    -> "Auf Nico aufgepasst und bewertet"
#    (keep the OLE_LINK4 bookmark around "aufgepasst" intact)
# ---------------------------------------------------------------------
$rng1 = $d.Content.Duplicate
$found1 = $rng1.Find.Execute(
    "Nico (Schnupperer) ", $false, $false, $false, $false, $false,
    $true, 1, $false, "Nico ", 2)
Write-Output "edit1: $found1"

# ---------------------------------------------------------------------
# 2) Table cell: "Meine" + " App " + "fertig programmiert" (3 runs)
#    -> a single run "Meine App fertig programmiert"
# ---------------------------------------------------------------------
$rng2 = $d.Content.Duplicate
$found2 = $rng2.Find.Execute(
    "Meine App fertig programmiert", $false, $false, $false, $false, $false,
    $true, 1, $false, "Meine App fertig programmiert", 2)
Write-Output "edit2: $found2"

# ---------------------------------------------------------------------
# 3) Wochenrueckblick paragraph: drop "Ein " before the mention of the
#    "Schnupperer" and rename him "Nico" instead.
# ---------------------------------------------------------------------
$rng3a = $d.Content.Duplicate
$t3a = "Diese Woche war eine sehr entspannte Woche. Wir hatten unseren letzten Schnuppertag, welcher sehr gut verlaufen ist. Ein "
$r3a = "Diese Woche war eine sehr entspannte Woche. Wir hatten unseren letzten Schnuppertag, welcher sehr gut verlaufen ist. "
$found3a = $rng3a.Find.Execute($t3a, $false, $false, $false, $false, $false,
    $true, 1, $false, $r3a, 2)
Write-Output "edit3a: $found3a"

$rng3b = $d.Content.Duplicate
$found3b = $rng3b.Find.Execute(
    "Schnupperer", $false, $false, $false, $false, $false,
    $true, 1, $false, "Nico ", 2)
Write-Output "edit3b: $found3b"

$rng3c = $d.Content.Duplicate
$t3c = " hat mich sogar in seiner Bewertung erwähnt, was ich sehr geschätzt habe. Ich finde, dass wir uns als Gruppe verbessert haben. Bei dem ersten Schnuppertag war ich auch dabei und sehe einen sehr gro"
$r3c = "hat mich sogar in seiner Bewertung erwähnt, was ich sehr geschätzt habe. Ich finde, dass wir uns als Gruppe verbessert haben. Bei dem ersten Schnuppertag war ich auch dabei und sehe einen sehr gro"
$found3c = $rng3c.Find.Execute($t3c, $false, $false, $false, $false, $false,
    $true, 1, $false, $r3c, 2)
Write-Output "edit3c: $found3c"
